$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 28; this pushes the existing rows 28..123 down to 29..124
# and Excel automatically copies formatting (e.g. the date style on column D) to the new row.
$ws.Rows("28:28").Insert()

# Copy the unchanged columns from the row directly below (which now holds the data that used
# to be in row 28) into the new row 28, then overwrite the columns that actually have new values.
$ws.Range("A28").Value = $ws.Range("A29").Value2
$ws.Range("B28").Value = $ws.Range("B29").Value2
$ws.Range("C28").Value = $ws.Range("C29").Value2
$ws.Range("E28").Value = $ws.Range("E29").Value2
$ws.Range("F28").Value = $ws.Range("F29").Value2
$ws.Range("G28").Value = $ws.Range("G29").Value2
$ws.Range("H28").Value = $ws.Range("H29").Value2
$ws.Range("I28").Value = $ws.Range("I29").Value2
$ws.Range("J28").Value = $ws.Range("J29").Value2
$ws.Range("N28").Value = $ws.Range("N29").Value2
$ws.Range("O28").Value = $ws.Range("O29").Value2
$ws.Range("Q28").Value = $ws.Range("Q29").Value2
$ws.Range("R28").Value = $ws.Range("R29").Value2

# New data for row 28
$ws.Range("D28").Value = 45145
$ws.Range("K28").Value = 3400
$ws.Range("L28").Value = 3500
$ws.Range("M28").Value = 3450
$ws.Range("P28").Value = 1725
